# Insert a new weekly record at row 472 (Hortaliza, Mercado Mayorista Lo
# Valledor de Santiago - Haba). Inserting the row pushes the former rows
# 472-504 down to 473-505, and the new row 472 receives the latest week's
# price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 472, shifting rows 472:504 -> 473:505
$ws.Rows("472:472").Insert()

# Populate the newly inserted row 472 with the new weekly observation
$ws.Range("A472").Value2 = 6
$ws.Range("B472").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C472").Value2 = "Metropolitana"
$ws.Range("D472").Value2 = 45265
$ws.Range("D472").NumberFormat = $ws.Range("D473").NumberFormat
$ws.Range("E472").Value2 = 13
$ws.Range("F472").Value2 = 100112026
$ws.Range("G472").Value2 = "Haba"
$ws.Range("H472").Value2 = "Sin especificar"
$ws.Range("I472").Value2 = "Primera"
$ws.Range("J472").Value2 = 230
$ws.Range("K472").Value2 = 12000
$ws.Range("L472").Value2 = 12000
$ws.Range("M472").Value2 = 12000
$ws.Range("N472").Value2 = "$/saco 25 kilos"
$ws.Range("O472").Value2 = "Región Metropolitana"
$ws.Range("P472").Value2 = 480
$ws.Range("Q472").Value2 = 25
$ws.Range("R472").Value2 = "Hortaliza"
